$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.741.04'
$ws.Range("E2").Value = '  -1.42%  '
$ws.Range("D3").Value = '2.542.36'
$ws.Range("E3").Value = '  -1.96%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '310.15'
$ws.Range("E5").Value = '  -2.05%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '101.53'
$ws.Range("E6").Value = '  +3.87%  '
$ws.Range("E7").Value = '  -1.15%  '
$ws.Range("E8").Value = '  +0.17%  '
$ws.Range("E9").Value = '  -2.35%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '36.39'
$ws.Range("E10").Value = '  +1.25%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0806'
$ws.Range("E11").Value = '  -1.29%  '
$ws.Range("E12").Value = '  -2.46%  '
$ws.Range("E13").Value = '  +0.27%  '
$ws.Range("D14").Value = '2.930.31'
$ws.Range("E14").Value = '  -2.03%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.76'
$ws.Range("E15").Value = '  +3.10%  '
$ws.Range("D16").Value = '2.559.93'
$ws.Range("E16").Value = '  -1.01%  '
$ws.Range("E17").Value = '  -4.18%  '
$ws.Range("D18").Value = '42.729.78'
$ws.Range("E18").Value = '  -1.57%  '
$ws.Range("E19").Value = '  -1.62%  '
$ws.Range("E20").Value = '  -1.31%  '
$ws.Range("E21").Value = '  -2.60%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '69.58'
$ws.Range("E22").Value = '  -0.34%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '244.86'
$ws.Range("E23").Value = '  -4.26%  '
$ws.Range("E24").Value = '  -2.50%  '
$ws.Range("E25").Value = '  -1.51%  '
$ws.Range("E26").Value = '  +0.02%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '26.19'
$ws.Range("E27").Value = '  -4.68%  '
$ws.Range("E28").Value = '  -4.80%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '39.44'
$ws.Range("E29").Value = '  -2.38%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '10.22'
$ws.Range("E30").Value = '  -1.31%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.81'
$ws.Range("E31").Value = '  -1.53%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '156.38'
$ws.Range("E32").Value = '  -0.78%  '
$ws.Range("E33").Value = '  +10.16%  '
$ws.Range("E34").Value = '  -2.15%  '
$ws.Range("E35").Value = '  -2.59%  '
$ws.Range("E36").Value = '  -5.36%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.19'
$ws.Range("E37").Value = '  -7.60%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '18.29'
$ws.Range("E38").Value = '  -2.30%  '
$ws.Range("E39").Value = '  -0.05%  '
$ws.Range("E40").Value = '  +0.09%  '
$ws.Range("E41").Value = '  +7.20%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '22.06'
$ws.Range("E42").Value = '  -1.90%  '
$ws.Range("E43").Value = '  +1.73%  '
$ws.Range("E44").Value = '  +0.08%  '
$ws.Range("E45").Value = '  -1.70%  '
$ws.Range("D46").Value = '1.985.51'
$ws.Range("E46").Value = '  -1.71%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.92'
$ws.Range("E47").Value = '  -0.98%  '
$ws.Range("D48").Value = '2.783.36'
$ws.Range("E48").Value = '  -1.94%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '80.90'
$ws.Range("E49").Value = '  -3.59%  '
$ws.Range("E50").Value = '  -0.82%  '
$ws.Range("B51").Value = 'SEI'
$ws.Range("C51").Value = 'https://coinranking.com/coin/8nxCqs-uj+sei-sei'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.850'
$ws.Range("E51").Value = '  +7.95%  '
